# Auto-generated edit script for keyframes.xlsx
# Commit message: '360 flip late kickflip'
#
# The keyframe table (row 1 = keyframe times, rows 2-28 = one animated
# channel per row) gains 11 new keyframe columns (AH:AR) appended after
# the previous last column (AG), and a small number of pre-existing
# values within B:AG are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. The new header cells (row 1, AH1:AR1) are keyframe-time values
#        and must keep the same bold / bordered / centered style as the
#        rest of row 1, so copy formatting from AG1 across first.
$ws.Range("AG1").Copy() | Out-Null
$ws.Range("AH1:AR1").PasteSpecial(-4122) | Out-Null

# --- 2. Build the new AH:AR values (28 rows x 11 cols) as a jagged array
#        and blit it into the sheet in a single Range write.
$newRows = @(
  @(50.6, 52, 52.9, 53.2, 57.2, 58.2, 59.2, 61.5, 62, 70, 76.59999999999999),
  @(0.2, -0.8, -0.1, 0, 0.6, 0.6, 0.6, 0.6, 0.6, 0, 0),
  @(0.2, 0.2, 0.1, 0.1, 0.5, 0.5, 0.5, 0.3, 0.2, 0.2, 0),
  @(1.1, 0.9, 0.6, 1, 1, 1, 1, 1, 1, 1, -5),
  @(0.2, 0.2, 0, 0, 0, 0, 0, 0, 0, 0, 3),
  @(0.1, 0, 0.1, 0.1, 0.4, 0.4, 0.4, 0.3, 0.1, 0.1, 1),
  @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, -10),
  @(3, 3, 3, 0.5, 0.5, 0.5, 0.5, 0.5, 3, 0.1, 0.1),
  @(0, 0, 0, 0.4, 4.5, 5.5, 6.6, 13.2, 12.6, 12.6, 12.6),
  @(0, 0, 0, 0.5, 4.5, 5.7, 6.1, 6.3, 6.3, 6.3, 6.3),
  @(0, 0, 0, -0.8, -0.4, -0.4, -0.2, 0.3, 0, 0, 0),
  @(0.2, 0.2, 0, 0, 0, 0, 0, 0, 0, 2, 2),
  @(0, 0, 0, 0.13, 0.3, 0.3, 0.2, 0.1, 0, 0, 0),
  @(0.1, 0.1, 0, 0, 0, 0, 0, 0.1, 0.1, 0.1, 0.1),
  @(-0.1, -0.1, 0, -0.1, -0.1, -0.1, -0.1, -0.1, -0.1, -0.1, -0.1),
  @(0.2, 0.2, -0.1, -0.1, -0.1, -0.1, -0.1, 0, 0, 2, 2),
  @(-0.02, -0.02, -0.31, 0, 0.3, 0.3, 0.3, 0.1, -0.02, -0.02, -0.02),
  @(0, 0, -0.2, 0, 0, 0, 0, 0, 0, 0, 0),
  @(0.7, 0.7, 0, 0, 0.5, 0.6, 0.6, 1.1, 0.7, 0.7, 0.7),
  @(0.1, 0, 0.43, -0.1, 0.3, 0.3, 0.2, 0.2, 0.1, 0.1, 0.1),
  @(0.3, 0.3, 0.2, 0.2, 0.2, 0.2, 0.2, 0.2, 0.2, 0.2, 0.2),
  @(0.2, 0, 0.8, 0.1, 0.6, 0.7, 0.8, 0.6, 0.2, 0.2, 0.2),
  @(0.1, -0.1, 1, -1, 0, 0, 0, 0.4, 0.2, 0.2, 0.2),
  @(-0.5, -0.5, 0.1, 0, 0, 0, 0, 0, -0.5, -0.5, -0.5),
  @(0.1, 0.1, 0.4, 0.5, 0.4, 0.2, -0.1, 0.2, 0.1, 0.1, 0.1),
  @(-0.1, -0.1, 0.1, 0, 0.4, 0.3, 0.3, 0.1, 0.1, 0.1, 0.1),
  @(0.1, 0.1, 0.8, 0.9, 0.1, 0.5, 0.1, 0.3, 0.1, 0.1, 0.1),
  @(-0.1, -0.1, 1.1, 0, 0, -0.3, -0.5, -0.8, -0.1, -0.1, -0.1)
)

$newBlock = New-Object 'object[,]' 28,11
for ($r = 0; $r -lt 28; $r++) {
  for ($c = 0; $c -lt 11; $c++) {
    $newBlock[$r,$c] = $newRows[$r][$c]
  }
}
$ws.Range("AH1:AR28").Value2 = $newBlock

# --- 3. A handful of pre-existing cells (within B:AG) changed value too.
$ws.Range("AE2").Value2 = 0.4
$ws.Range("AE3").Value2 = 0.2
$ws.Range("AE4").Value2 = 1.1
$ws.Range("N6").Value2 = 0.3
$ws.Range("N15").Value2 = -0.1
$ws.Range("N17").Value2 = 0
$ws.Range("N19").Value2 = -0.3
$ws.Range("O19").Value2 = 0.3
$ws.Range("N20").Value2 = -0.1
$ws.Range("O20").Value2 = 0.1
$ws.Range("N22").Value2 = 0.1
$ws.Range("O22").Value2 = 0.6
$ws.Range("N23").Value2 = -1
$ws.Range("N28").Value2 = 0
